$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.442.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.50%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.988.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.40%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -9.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.589'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.29%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.82'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.371'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.42'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.51%  '

$ws.Range("E11").Value = '  -3.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0989'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.92%  '

$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.282.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.28%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.11'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.61%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.41'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.758'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.92%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.988.67'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.10%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.426.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '67.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.13%  '

$ws.Range("D21").Value = '0.0₃0804'
$ws.Range("E21").Value = '  -4.27%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.25'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '221.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.11%  '

$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("E26").Value = '  -9.84%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.67'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.00%  '

$ws.Range("E29").Value = '  -1.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.87'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.33'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.117'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0605'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.26'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.34'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.16%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.43'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.14%  '

$ws.Range("B38").Value = 'BinanceUSD'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.03%  '

$ws.Range("E39").Value = '  -3.35%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.60'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.36%  '

$ws.Range("E41").Value = '  -1.10%  '

$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0941'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.61%  '

$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.453.82'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0203'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.13'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.991'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.62%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.88'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.85'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.68'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.79%  '
